{"js": "// Update the date line (first paragraph in the body, outside the table).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst dateParagraph = paragraphs.items[0];\ndateParagraph.getRange().insertText(\"2024-06-01 Saturday\", Word.InsertLocation.replace);\n\n// Update every equation cell in the single 20x5 table, in row-major order.\nconst newValues = [\"57+24=81\", \"33+42=75\", \"92-49=43\", \"24-7=17\", \"55-29=26\", \"38+24=62\", \"19+71=90\", \"40+13=53\", \"33+43=76\", \"56+9=65\", \"43-14=29\", \"85-82=3\", \"59+29=88\", \"4+60=64\", \"91-4=87\", \"79+14=93\", \"11+37=48\", \"43+50=93\", \"2+38=40\", \"72-40=32\", \"1+20=21\", \"54-18=36\", \"51+31=82\", \"97-76=21\", \"92-3=89\", \"63-8=55\", \"43+44=87\", \"90-38=52\", \"22+37=59\", \"54-51=3\", \"74-61=13\", \"17+47=64\", \"38+3=41\", \"26+68=94\", \"41+57=98\", \"75-68=7\", \"5-3=2\", \"33+57=90\", \"9+63=72\", \"70-50=20\", \"93-22=71\", \"49+37=86\", \"19+30=49\", \"56-44=12\", \"14+68=82\", \"92-53=39\", \"45-17=28\", \"88-20=68\", \"49-17=32\", \"47+45=92\", \"57-49=8\", \"85-24=61\", \"44-25=19\", \"63-43=20\", \"47-39=8\", \"74+11=85\", \"76+19=95\", \"19+23=42\", \"97-56=41\", \"52+3=55\", \"61-40=21\", \"2+58=60\", \"69-11=58\", \"70-21=49\", \"1+32=33\", \"23+8=31\", \"51+4=55\", \"57-20=37\", \"13+31=44\", \"5+8=13\", \"81-47=34\", \"17+9=26\", \"61-35=26\", \"79-65=14\", \"63+12=75\", \"56-2=54\", \"72-31=41\", \"3+11=14\", \"23+48=71\", \"65+11=76\", \"85-49=36\", \"82-38=44\", \"52+38=90\", \"4+90=94\", \"45-1=44\", \"25+12=37\", \"52-51=1\", \"65+4=69\", \"48+23=71\", \"35+25=60\", \"14+36=50\", \"92-50=42\", \"84-62=22\", \"32+40=72\", \"83+14=97\", \"55+5=60\", \"7+45=52\", \"67+14=81\", \"26+42=68\", \"76-29=47\"];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\nconst columnCount = 5;\nlet index = 0;\nfor (let row = 0; row < table.rowCount; row++) {\n  for (let col = 0; col < columnCount; col++) {\n    const cell = table.getCell(row, col);\n    const cellParagraphs = cell.body.paragraphs;\n    cellParagraphs.load(\"items\");\n    await context.sync();\n    const cellParagraph = cellParagraphs.items[0];\n    cellParagraph.getRange().insertText(newValues[index], Word.InsertLocation.replace);\n    index++;\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Update the date line (first paragraph in the body, outside the table).\n$d.Paragraphs.Item(1).Range.Text = \"2024-06-01 Saturday\"\n\n# Update every equation cell in the single 20x5 table, in row-major order.\n$newValues = @(\"57+24=81\",\"33+42=75\",\"92-49=43\",\"24-7=17\",\"55-29=26\",\"38+24=62\",\"19+71=90\",\"40+13=53\",\"33+43=76\",\"56+9=65\",\"43-14=29\",\"85-82=3\",\"59+29=88\",\"4+60=64\",\"91-4=87\",\"79+14=93\",\"11+37=48\",\"43+50=93\",\"2+38=40\",\"72-40=32\",\"1+20=21\",\"54-18=36\",\"51+31=82\",\"97-76=21\",\"92-3=89\",\"63-8=55\",\"43+44=87\",\"90-38=52\",\"22+37=59\",\"54-51=3\",\"74-61=13\",\"17+47=64\",\"38+3=41\",\"26+68=94\",\"41+57=98\",\"75-68=7\",\"5-3=2\",\"33+57=90\",\"9+63=72\",\"70-50=20\",\"93-22=71\",\"49+37=86\",\"19+30=49\",\"56-44=12\",\"14+68=82\",\"92-53=39\",\"45-17=28\",\"88-20=68\",\"49-17=32\",\"47+45=92\",\"57-49=8\",\"85-24=61\",\"44-25=19\",\"63-43=20\",\"47-39=8\",\"74+11=85\",\"76+19=95\",\"19+23=42\",\"97-56=41\",\"52+3=55\",\"61-40=21\",\"2+58=60\",\"69-11=58\",\"70-21=49\",\"1+32=33\",\"23+8=31\",\"51+4=55\",\"57-20=37\",\"13+31=44\",\"5+8=13\",\"81-47=34\",\"17+9=26\",\"61-35=26\",\"79-65=14\",\"63+12=75\",\"56-2=54\",\"72-31=41\",\"3+11=14\",\"23+48=71\",\"65+11=76\",\"85-49=36\",\"82-38=44\",\"52+38=90\",\"4+90=94\",\"45-1=44\",\"25+12=37\",\"52-51=1\",\"65+4=69\",\"48+23=71\",\"35+25=60\",\"14+36=50\",\"92-50=42\",\"84-62=22\",\"32+40=72\",\"83+14=97\",\"55+5=60\",\"7+45=52\",\"67+14=81\",\"26+42=68\",\"76-29=47\")\n\n$table = $d.Tables.Item(1)\n$columnCount = 5\n$index = 0\nfor ($row = 1; $row -le $table.Rows.Count; $row++) {\n  for ($col = 1; $col -le $columnCount; $col++) {\n    $table.Cell($row, $col).Range.Text = $newValues[$index]\n    $index++\n  }\n}\n"}
